# edit.ps1
# Applies the "fix: ready to present" change set:
#  - LOW_PRICE (ANGKOST):  add row 3 (new low-price pick)
#  - HIGH_PRICE (ANGKOST): replace row 2 with the new top pick, add row 3
#  - LOW_SCORE (ANGKOST):  replace row 2 with the new top pick, add row 3
#  - HIGH_SCORE (ANGKOST): bump row 2's score (F2), add row 3
#  - add a brand-new RECOMMENDATION (ANGKOST) sheet with the two
#    highest-recommendation-score products

function Set-TextCell {
    # Force a cell to be written as TEXT even when the string looks
    # numeric (e.g. "4.6"), then strip the number-format override so the
    # cell is left with no extra style applied.
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

function Set-ProductRow {
    param($ws, $rowNum, $id, $name, $price, $rating, $link, $score)
    $ws.Cells.Item($rowNum, 1).Value = $id
    $ws.Cells.Item($rowNum, 2).Value = $name
    # NB: compute the addresses into locals first -- passing an inline
    # ("C" + $rowNum) expression straight through as a nested call's
    # argument silently drops the call in this interpreter.
    $cAddr = "C" + $rowNum
    $dAddr = "D" + $rowNum
    Set-TextCell $ws $cAddr $price
    Set-TextCell $ws $dAddr $rating
    $ws.Cells.Item($rowNum, 5).Value = $link
    $ws.Cells.Item($rowNum, 6).Value = $score
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# LOW_PRICE (ANGKOST): add row 3
# ---------------------------------------------------------------------
$wsLowPrice = $wb.Worksheets.Item("LOW_PRICE (ANGKOST)")
Set-ProductRow $wsLowPrice 3 1414556557 `
    "Blinqshop Baju Tidur Piyama Motid Aurel" `
    "Rp27.900" "4.6" `
    "https://www.tokopedia.com/blinqshop/blinqshop-baju-tidur-piyama-motid-aurel-dst-putih?extParam=ivf%3Dfalse%26src%3Dsearch" `
    1048.8

# ---------------------------------------------------------------------
# HIGH_PRICE (ANGKOST): replace row 2, add row 3
# ---------------------------------------------------------------------
$wsHighPrice = $wb.Worksheets.Item("HIGH_PRICE (ANGKOST)")
Set-ProductRow $wsHighPrice 2 1862123531 `
    "Logitech G304 Lightspeed Wireless Gaming Mouse - Logitech G-304" `
    "Rp498.000" "5.0" `
    "https://www.tokopedia.com/duniacom-srv/logitech-g304-lightspeed-wireless-gaming-mouse-logitech-g-304-putih?extParam=ivf%3Dfalse%26src%3Dsearch" `
    4055
Set-ProductRow $wsHighPrice 3 1722207942 `
    "Baju Tidur Wanita Setelan Panjang / Piyama Set Kimono" `
    "Rp118.340" "4.7" `
    "https://www.tokopedia.com/twinolshop-2/baju-tidur-wanita-setelan-panjang-piyama-set-kimono-pink-leaf?extParam=ivf%3Dfalse%26src%3Dsearch" `
    4107.8

# ---------------------------------------------------------------------
# LOW_SCORE (ANGKOST): replace row 2, add row 3
# ---------------------------------------------------------------------
$wsLowScore = $wb.Worksheets.Item("LOW_SCORE (ANGKOST)")
Set-ProductRow $wsLowScore 2 1862123531 `
    "Logitech G304 Lightspeed Wireless Gaming Mouse - Logitech G-304" `
    "Rp498.000" "5.0" `
    "https://www.tokopedia.com/duniacom-srv/logitech-g304-lightspeed-wireless-gaming-mouse-logitech-g-304-putih?extParam=ivf%3Dfalse%26src%3Dsearch" `
    4055
Set-ProductRow $wsLowScore 3 1977206481 `
    "Lingerie Sexy Outer+Dress Transparan Piyama Baju Tidur Kimono Lr11" `
    "Rp92.000" "4.9" `
    "https://www.tokopedia.com/officialbianglalaid/lingerie-sexy-outer-dress-transparan-piyama-baju-tidur-kimono-lr11-merah?extParam=ivf%3Dfalse%26src%3Dsearch" `
    955.5000000000001

# ---------------------------------------------------------------------
# HIGH_SCORE (ANGKOST): bump row 2's score only, add row 3
# ---------------------------------------------------------------------
$wsHighScore = $wb.Worksheets.Item("HIGH_SCORE (ANGKOST)")
$wsHighScore.Cells.Item(2, 6).Value = 173626.6
Set-ProductRow $wsHighScore 3 1722207942 `
    "Baju Tidur Wanita Setelan Panjang / Piyama Set Kimono" `
    "Rp118.340" "4.7" `
    "https://www.tokopedia.com/twinolshop-2/baju-tidur-wanita-setelan-panjang-piyama-set-kimono-pink-leaf?extParam=ivf%3Dfalse%26src%3Dsearch" `
    4107.8

# ---------------------------------------------------------------------
# New sheet: RECOMMENDATION (ANGKOST), appended after the last sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRec = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsRec.Name = "RECOMMENDATION (ANGKOST)"
$wsRec.Tab.Color = 50431   # 0x00C4FF -> OLE BGR for RGB(FF,C4,00)

# Match the sheetPr / page-setup conventions used by the other sheets.
$wsRec.Outline.SummaryRow = 1
$wsRec.Outline.SummaryColumn = 1
$wsRec.PageSetup.LeftMargin = 54
$wsRec.PageSetup.RightMargin = 54
$wsRec.PageSetup.TopMargin = 72
$wsRec.PageSetup.BottomMargin = 72
$wsRec.PageSetup.HeaderMargin = 36
$wsRec.PageSetup.FooterMargin = 36

$wsRec.Cells.Item(1, 1).Value = "Id Produk"
$wsRec.Cells.Item(1, 2).Value = "Nama Produk"
$wsRec.Cells.Item(1, 3).Value = "Harga Produk"
$wsRec.Cells.Item(1, 4).Value = "Bintang Produk"
$wsRec.Cells.Item(1, 5).Value = "Link Produk"
$wsRec.Cells.Item(1, 6).Value = "Skor Produk (Rekomendasi)"

Set-ProductRow $wsRec 2 1862123531 `
    "Logitech G304 Lightspeed Wireless Gaming Mouse - Logitech G-304" `
    "Rp498.000" "5.0" `
    "https://www.tokopedia.com/duniacom-srv/logitech-g304-lightspeed-wireless-gaming-mouse-logitech-g-304-putih?extParam=ivf%3Dfalse%26src%3Dsearch" `
    6.973817219257057
Set-ProductRow $wsRec 3 1977206481 `
    "Lingerie Sexy Outer+Dress Transparan Piyama Baju Tidur Kimono Lr11" `
    "Rp92.000" "4.9" `
    "https://www.tokopedia.com/officialbianglalaid/lingerie-sexy-outer-dress-transparan-piyama-baju-tidur-kimono-lr11-merah?extParam=ivf%3Dfalse%26src%3Dsearch" `
    3.482186883868269

$wsLowPrice.Range("A1").Select()
